# Fruta / hortaliza, semanal
# Insert two new weekly records (Primera/Segunda, Paine, 15kg crates) ahead
# of the existing "Vega Central Mapocho de Santiago - Granada" rows, pushing
# the previously-last "Segunda" (O'Higgins, 18kg) record into its own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 31-36 down to 33-38, leaving two blank rows (31:32) for the
# new data.
$ws.Rows("31:32").Insert()

# New row 31: Primera / Paine / 15 kilos
$ws.Range("A31").Value = 9
$ws.Range("B31").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C31").Value = "Metropolitana"
$ws.Range("D31").Value = 45034
$ws.Range("E31").Value = 13
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100104
$ws.Range("H31").Value = "Frutos de pepita"
$ws.Range("I31").Value = 100104001
$ws.Range("J31").Value = "Granada"
$ws.Range("K31").Value = "Wonderfull"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 300
$ws.Range("N31").Value = 9000
$ws.Range("O31").Value = 9000
$ws.Range("P31").Value = 9000
$ws.Range("Q31").Value = "$/caja 15 kilos granel"
$ws.Range("R31").Value = "Paine"
$ws.Range("S31").Value = 600
$ws.Range("T31").Value = 15

# New row 32: Segunda / Paine / 15 kilos
$ws.Range("A32").Value = 9
$ws.Range("B32").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 45034
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100104
$ws.Range("H32").Value = "Frutos de pepita"
$ws.Range("I32").Value = 100104001
$ws.Range("J32").Value = "Granada"
$ws.Range("K32").Value = "Wonderfull"
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 280
$ws.Range("N32").Value = 6000
$ws.Range("O32").Value = 6000
$ws.Range("P32").Value = 6000
$ws.Range("Q32").Value = "$/caja 15 kilos granel"
$ws.Range("R32").Value = "Paine"
$ws.Range("S32").Value = 400
$ws.Range("T32").Value = 15
